$wb = $excel.ActiveWorkbook

# --- NEW-PLAN sheet: update the plan row and reshape the small table ---
$ws5 = $wb.Worksheets.Item("NEW-PLAN")

# Row 2 gets a new plan code / MSISDN / ICCID-like values.
# (E2 and D2 are written in this order so any newly-introduced shared
# strings land in the same order as the authored edit.)
$ws5.Cells.Item(2, 3).Value = "PLC"
$ws5.Cells.Item(2, 5).Value = "94754900"
$ws5.Cells.Item(2, 4).Value = "89598076101039725728"

# The old row 3 (uat / Puper3 / PLTTEF / ...) is removed entirely.
$ws5.Rows.Item(3).Delete()

# Three blank rows are added back underneath, keeping the D/E column
# formatting from row 2 (so D3:D5/E3:E5 stay styled but empty).
$ws5.Range("D2:E2").Copy()
$ws5.Range("D3:E3").PasteSpecial(-4122)
$ws5.Range("D4:E4").PasteSpecial(-4122)
$ws5.Range("D5:E5").PasteSpecial(-4122)

# Selection on NEW-PLAN moves to E8, and the sheet is no longer the
# selected tab.
$ws5.Range("E8").Select() | Out-Null

# --- Make POTENTIAL-RESCLIENTS the active/selected worksheet tab ---
$ws1 = $wb.Worksheets.Item("POTENTIAL-RESCLIENTS")
$ws1.Select() | Out-Null
